$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.703.90"
$ws.Range("E2").Value = "  -0.60%  "

# Row 3
$ws.Range("D3").Value = "2.216.69"
$ws.Range("E3").Value = "  -0.48%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").Value = "'252.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.01%  "

# Row 6
$ws.Range("D6").Value = "'0.626"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "

# Row 7
$ws.Range("D7").Value = "'70.39"
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = "  -0.14%  "

# Row 9
$ws.Range("D9").Value = "'0.598"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.51%  "

# Row 10
$ws.Range("D10").Value = "'39.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.39%  "

# Row 11
$ws.Range("D11").Value = "'0.0967"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.70%  "

# Row 12
$ws.Range("D12").Value = "'58.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.32%  "

# Row 13
$ws.Range("D13").Value = "'7.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.48%  "

# Row 14
$ws.Range("E14").Value = "  -0.10%  "

# Row 15
$ws.Range("D15").Value = "2.547.26"
$ws.Range("E15").Value = "  -0.75%  "

# Row 16
$ws.Range("D16").Value = "'15.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.60%  "

# Row 17
$ws.Range("D17").Value = "'0.896"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.17%  "

# Row 18
$ws.Range("D18").Value = "2.225.73"
$ws.Range("E18").Value = "  +0.18%  "

# Row 19
$ws.Range("D19").Value = "41.716.13"
$ws.Range("E19").Value = "  -0.38%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0966"
$ws.Range("E20").Value = "  +0.91%  "

# Row 21
$ws.Range("E21").Value = "  +1.98%  "

# Row 22
$ws.Range("D22").Value = "'72.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.01%  "

# Row 23
$ws.Range("E23").Value = "  +0.15%  "

# Row 24
$ws.Range("E24").Value = "  +1.96%  "

# Row 25
$ws.Range("E25").Value = "  +13.24%  "

# Row 26
$ws.Range("D26").Value = "'12.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +23.26%  "

# Row 27
$ws.Range("E27").Value = "  -0.01%  "

# Row 28
$ws.Range("E28").Value = "  +3.87%  "

# Row 29
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'171.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.30%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.84%  "

# Row 31
$ws.Range("D31").Value = "'20.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.38%  "

# Row 32
$ws.Range("E32").Value = "  +2.53%  "

# Row 33
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'5.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.47%  "

# Row 34
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "'0.124"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.80%  "

# Row 35
$ws.Range("D35").Value = "'0.0745"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.73%  "

# Row 36
$ws.Range("D36").Value = "'4.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.44%  "

# Row 37
$ws.Range("D37").Value = "'25.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.46%  "

# Row 38
$ws.Range("D38").Value = "'4.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.81%  "

# Row 39
$ws.Range("D39").Value = "'0.0307"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.54%  "

# Row 40
$ws.Range("E40").Value = "  -0.31%  "

# Row 41
$ws.Range("D41").Value = "'5.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.10%  "

# Row 42
$ws.Range("D42").Value = "'12.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +26.80%  "

# Row 43
$ws.Range("D43").Value = "'65.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.38%  "

# Row 44
$ws.Range("D44").Value = "'0.205"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.44%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'8.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.61%  "

# Row 46
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").Value = "'4.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.85%  "

# Row 47
$ws.Range("D47").Value = "'0.102"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.65%  "

# Row 48
$ws.Range("E48").Value = "  +0.26%  "

# Row 49
$ws.Range("D49").Value = "'4.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.95%  "

# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'2.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.41%  "

# Row 51
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "'1.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.84%  "
